$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 70: 四方坪站 (site 2 / shared string index 2)
$ws.Cells.Item(70, 1).Value = 45935
$ws.Cells.Item(70, 2).Value = "四方坪站"
$ws.Cells.Item(70, 3).Formula = "=15681/127"
$ws.Cells.Item(70, 4).Formula = "=C70/(24*60)"
$ws.Cells.Item(70, 5).Formula = "=9133.79/127"
$ws.Cells.Item(70, 6).Formula = "=3141.45/127"
$ws.Cells.Item(70, 7).Formula = "=9133.79/(15681/60)"
$ws.Cells.Item(70, 8).Formula = "=385/127"

# Row 71: 高岭站 (site 3 / shared string index 3)
$ws.Cells.Item(71, 1).Value = 45935
$ws.Cells.Item(71, 2).Value = "高岭站"
$ws.Cells.Item(71, 3).Formula = "=4430/36"
$ws.Cells.Item(71, 4).Formula = "=C71/(24*60)"
$ws.Cells.Item(71, 5).Formula = "=3646.37/36"
$ws.Cells.Item(71, 6).Formula = "=912.3/36"
$ws.Cells.Item(71, 7).Formula = "=3646.37/(4430/60)"
$ws.Cells.Item(71, 8).Formula = "=114/36"

# Update the active selection to match the new end of the data (H73)
$ws.Range("H73").Select()
